$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- "category" column, inserted right after "property_category" (col H) ---
# This pushes the former I/J/K ("date"/"legislator_name"/"legislator_id")
# one column to the right, to J/K/L.
$ws.Columns.Item(9).Insert()
$ws.Range("I1").Value = "category"
$ws.Range("I2").Value = "normal"

# --- "source_file" / "index" columns, appended after "legislator_id" (col L) ---
$ws.Columns.Item(13).Insert()
$ws.Columns.Item(14).Insert()

$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

$ws.Range("M2").Value = "tmpc6841"
$ws.Range("N2").Value = 74
